$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Merge the split runs "insert" + "_sort" -> "insert_sort" and the three
#    trailing runs describing the insert_sort function into a single run.
#    Scope the Find to the specific paragraph so we don't touch the
#    unrelated, already-combined "insert_sort" text earlier in the doc
#    (which stays split into "insert_s" + "ort" there).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like 'The "insert_sort" function accepts*') {
        $rng = $p.Range
        $rng.Find.Execute("insert_sort", $false, $false, $false, $false, $false, $true, 1, $false, "insert_sort", 2) | Out-Null
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Remove the old "_GoBack" bookmark that wraps "LinkList1" near the top
#    of the document.
# ---------------------------------------------------------------------------
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

# ---------------------------------------------------------------------------
# 3. Delete the large trailing block of paragraphs describing the "Create a
#    Linked List ... library" task, the "Bounce-Task" section, and the
#    bounce-mark note, leaving only the final "Note:" paragraph's shell
#    (its pPr/rPr), which becomes empty.
# ---------------------------------------------------------------------------
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like ' Create a Linked *') {
        $endPara = $i
    }
}
# Find the empty paragraph immediately preceding the "Create a Linked" one;
# there are two consecutive empty paragraphs there, and only the second one
# (immediately before "Create a Linked") is removed along with everything
# through the "Note:" paragraph's content.
$createIdx = $endPara
$startIdx = $createIdx - 1

$rngStart = $d.Paragraphs.Item($startIdx).Range.Start
# last paragraph (the "Note:" one) - its mark start is the end of the block.
$noteIdx = $d.Paragraphs.Count
$rngEnd = $d.Paragraphs.Item($noteIdx).Range.Start

$bigRange = $d.Range($rngStart, $rngEnd)
$bigRange.Delete()

# ---------------------------------------------------------------------------
# 4. Re-create the "_GoBack" bookmark collapsed at the start of the now
#    trailing "Note:" paragraph, then clear that paragraph's text, leaving
#    only its paragraph mark (pPr with <w:i/>) and the bookmark.
#    The bookmark must be added *before* the text is removed, because this
#    runtime's Bookmarks.Add mis-places bookmarks added into a completely
#    run-less (empty) paragraph.
# ---------------------------------------------------------------------------
$lastIdx = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIdx)
$bmPos = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmPos) | Out-Null

$lastPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$textRange = $d.Range($lastPara2.Range.Start, $lastPara2.Range.End - 1)
if ($textRange.Start -lt $textRange.End) {
    $textRange.Delete()
}
